# Insert a new data row at row 143 (pushing the existing rows 143-208 down
# to 144-209) and fill it in with a new price record for "Sandia" (Hortaliza),
# matching the commit "Fruta / hortaliza, semanal" weekly data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 143; existing rows shift down by one
# (Excel default behaviour for EntireRow.Insert is ShiftDown).
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new record.
$ws.Range("A143").Value = 9
$ws.Range("B143").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C143").Value = "Metropolitana"
$ws.Range("D143").Value = 44488
$ws.Range("E143").Value = 13
$ws.Range("F143").Value = 100112028
$ws.Range("G143").Value = "Sandia"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 250
$ws.Range("K143").Value = 800
$ws.Range("L143").Value = 1000
$ws.Range("M143").Value = 900
$ws.Range("N143").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O143").Value = "Perú"
$ws.Range("P143").Value = 900
$ws.Range("Q143").Value = 1
$ws.Range("R143").Value = "Hortaliza"
